# Remove the three subject rows that were dropped from the summary sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$namesToRemove = @(
    "COMPUTO FLEXIBLE (SOFTCOMPUTING)",
    "PROYECTO DE GESTION DE LA TECNOLOGIA DE INFORMACION",
    "PROYECTO DE SISTEMAS ROBUSTOS, PARALELOS Y DISTRIBUIDOS"
)

# Find the last used row in column A.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

foreach ($name in $namesToRemove) {
    for ($r = $lastRow; $r -ge 2; $r--) {
        $cellValue = $ws.Cells.Item($r, 1).Value2
        if ($cellValue -eq $name) {
            $ws.Rows.Item($r).Delete()
            break
        }
    }
}
